$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers I1 ("I0") and J1 ("IF"), matching the style of H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate I2:J54 with the new data
$data = New-Object 'object[,]' 53,2
$data[0,0] = 8
$data[0,1] = 8
$data[1,0] = 8
$data[1,1] = 8
$data[2,0] = 7
$data[2,1] = 7
$data[3,0] = 7
$data[3,1] = 7
$data[4,0] = 6
$data[4,1] = 6
$data[5,0] = 7
$data[5,1] = 8
$data[6,0] = 8
$data[6,1] = 8
$data[7,0] = 8
$data[7,1] = 8
$data[8,0] = 6
$data[8,1] = 6
$data[9,0] = 9
$data[9,1] = 9
$data[10,0] = 5
$data[10,1] = 5
$data[11,0] = 7
$data[11,1] = 9
$data[12,0] = 6
$data[12,1] = 7
$data[13,0] = 9
$data[13,1] = 9
$data[14,0] = 7
$data[14,1] = 7
$data[15,0] = 9
$data[15,1] = 9
$data[16,0] = 6
$data[16,1] = 6
$data[17,0] = 9
$data[17,1] = 9
$data[18,0] = 7
$data[18,1] = 7
$data[19,0] = 6
$data[19,1] = 6
$data[20,0] = 1
$data[20,1] = 3
$data[21,0] = 9
$data[21,1] = 9
$data[22,0] = 6
$data[22,1] = 7
$data[23,0] = 7
$data[23,1] = 7
$data[24,0] = 8
$data[24,1] = 8
$data[25,0] = 1
$data[25,1] = 2
$data[26,0] = 5
$data[26,1] = 6
$data[27,0] = 8
$data[27,1] = 8
$data[28,0] = 5
$data[28,1] = 6
$data[29,0] = 3
$data[29,1] = 6
$data[30,0] = 5
$data[30,1] = 7
$data[31,0] = 6
$data[31,1] = 6
$data[32,0] = 5
$data[32,1] = 6
$data[33,0] = 5
$data[33,1] = 5
$data[34,0] = 4
$data[34,1] = 6
$data[35,0] = 5
$data[35,1] = 7
$data[36,0] = 4
$data[36,1] = 6
$data[37,0] = 8
$data[37,1] = 9
$data[38,0] = 4
$data[38,1] = 5
$data[39,0] = 8
$data[39,1] = 9
$data[40,0] = 7
$data[40,1] = 8
$data[41,0] = 3
$data[41,1] = 5
$data[42,0] = 7
$data[42,1] = 8
$data[43,0] = 7
$data[43,1] = 8
$data[44,0] = 7
$data[44,1] = 7
$data[45,0] = 8
$data[45,1] = 8
$data[46,0] = 5
$data[46,1] = 6
$data[47,0] = 3
$data[47,1] = 7
$data[48,0] = 5
$data[48,1] = 5
$data[49,0] = 7
$data[49,1] = 8
$data[50,0] = 4
$data[50,1] = 5
$data[51,0] = 5
$data[51,1] = 6
$data[52,0] = 3
$data[52,1] = 3
$ws.Range("I2:J54").Value = $data

